$d = $word.ActiveDocument

function Set-ParaText($para, [string]$text) {
    $r = $para.Range
    $r.MoveEnd(1, -1) | Out-Null
    $r.Text = $text
}

# ------------------------------------------------------------------
# Work from the bottom of the document upwards so paragraph indices
# of not-yet-processed paragraphs are not disturbed by insertions.
# ------------------------------------------------------------------

# Paragraph 23 ("1-5a: ..." -> "1-6a: ...")
$p23 = $d.Paragraphs.Item(23)
$found = $p23.Range.Find.Execute("1-5a", $true, $false, $false, $false, $false, $true, 1, $false, "1-6a", 2)

# Paragraph 21 ("5a: L'utilisateur décide de ne pas confirmer son choix" -> "6a: ...")
$p21 = $d.Paragraphs.Item(21)
Set-ParaText $p21 "6a: L’utilisateur décide de ne pas confirmer son choix"

# Paragraph 20: tab + bold red "?" run -> two tabs + new message, then
# insert two brand-new paragraphs after it.
$p20 = $d.Paragraphs.Item(20)
Set-ParaText $p20 "`t`tLe système renvoie à l’utilisateur un message d’erreur et reste à cette étape sans prendre en compte la définition des plages horaires faites"

$p20 = $d.Paragraphs.Item(20)
$p20.Range.InsertParagraphAfter()
$p21new = $d.Paragraphs.Item(21)
Set-ParaText $p21new "5a: La modification des horaires d’arrivés provoque le non respect de la contrainte des plages horaires"
$p21new.Format.LeftIndent = 36
$p21new.Format.FirstLineIndent = 0

$p21new.Range.InsertParagraphAfter()
$p22new = $d.Paragraphs.Item(22)
Set-ParaText $p22new "Le système met en surbrillance les plages horaires non valide"
$p22new.Format.LeftIndent = 72
$p22new.Format.FirstLineIndent = 36

# Paragraph 19 ("4a. La nouvelle tournée calculée ne respecte pas ..." -> "4b: ...")
$p19 = $d.Paragraphs.Item(19)
Set-ParaText $p19 "4b: L’utilisateur définit une plage horaire qui n’est pas dans le système “horaire” "

# Paragraph 18 (existing tab + message -> two tabs + new message)
$p18 = $d.Paragraphs.Item(18)
Set-ParaText $p18 "`t`tLe système renvoie à l’utilisateur un message d’erreur et retourne à l’étape 3"

# Paragraph 17 ("2c. L'utilisateur donne un point de livraison ..." -> "4a: ...")
$p17 = $d.Paragraphs.Item(17)
Set-ParaText $p17 "4a: L’utilisateur sélectionne l’entrepôt comme livraison suivante"

# Paragraph 16 (existing tab + message -> two tabs + new message)
$p16 = $d.Paragraphs.Item(16)
Set-ParaText $p16 "`t`tLe système renvoie à l’utilisateur un message d’erreur et retourne à l’étape 1"

# Paragraph 15 ("2b. L'utilisateur donne un point de livraison qui existe déjà ..." -> "2b. L'utilisateur sélectionne ...")
$p15 = $d.Paragraphs.Item(15)
Set-ParaText $p15 "2b. L’utilisateur sélectionne un point sur la carte qui est déjà un point de livraison"

# Paragraph 14 (indent changes: remove left indent; add two tabs + new message)
$p14 = $d.Paragraphs.Item(14)
Set-ParaText $p14 "`t`tLe système renvoie à l’utilisateur un message d’erreur et retourne à l’étape 1"
$p14.Format.LeftIndent = 0
$p14.Format.FirstLineIndent = 36

# Paragraph 13 ("2a. L'utilisateur ne donne pas toutes les informations ..." -> "2a.L'utilisateur sélectionne ...")
$p13 = $d.Paragraphs.Item(13)
Set-ParaText $p13 "2a.L’utilisateur sélectionne un point sur la carte qui n’est pas un tronçon"

# Paragraph 10: new text, then insert a new numbered paragraph after it.
$p10 = $d.Paragraphs.Item(10)
Set-ParaText $p10 "Le système calcule le plus court chemin entre le nouveau point de livraison et la livraison avant et après le nouveau point de livraison puis met à jour les horaires d’arrivées"
$p10.Range.InsertParagraphAfter()
$p11new = $d.Paragraphs.Item(11)
Set-ParaText $p11new "L’utilisateur confirme le choix de la modification effectuée"

# Paragraph 9
$p9 = $d.Paragraphs.Item(9)
Set-ParaText $p9 "L’utilisateur via la carte ou le tableau sélectionne la livraison suivante puis définit les plages horaires de la nouvelle livraison"

# Paragraph 8
$p8 = $d.Paragraphs.Item(8)
Set-ParaText $p8 "Le système ajoute la nouvelle livraison sur la carte et demande à l’utilisateur avant quelle livraison elle se situe"

# Paragraph 7
$p7 = $d.Paragraphs.Item(7)
Set-ParaText $p7 "L’utilisateur choisit le tronçon sur la carte où il souhaite mettre la nouvelle livraison"

# Paragraph 6: merge the 3 runs (including the bold red "?") into a single run.
$p6 = $d.Paragraphs.Item(6)
Set-ParaText $p6 "Le système demande à l’utilisateur de placer la nouvelle livraison directement sur la carte en cliquant sur le tronçon voulu"

# Paragraph 3 (Precondition): remove "(non) " before "valide"
$d.Content.Find.Execute("(non) valide", $true, $false, $false, $false, $false, $true, 1, $false, "valide", 2) | Out-Null
